$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.703.17'
$ws.Range("E2").Value = '  +2.78%  '
$ws.Range("D3").Value = '1.864.64'
$ws.Range("E3").Value = '  +2.06%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.95%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7017'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.40%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07773'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3082'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.29%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.76'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.52%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07837'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.30%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.180'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.90%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.861.99'
$ws.Range("E13").Value = '  +1.82%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '92.98'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6963'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.642'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.62%  '
$ws.Range("D17").Value = '29.712.52'
$ws.Range("E17").Value = '  +2.79%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008400'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.69%  '
$ws.Range("D19").Value = '2.114.74'
$ws.Range("E19").Value = '  +1.53%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '244.02'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.82'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.657'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.69%  '
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1520'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.85%  '
$ws.Range("E26").Value = '  +3.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.37'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.22%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.40'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.72%  '
$ws.Range("E29").Value = '  +1.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.279'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.209'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.46%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.200'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.05%  '
$ws.Range("E33").Value = '  +0.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7894'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.919'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.162'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.99%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.692'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("D38").Value = '1.344.29'
$ws.Range("E38").Value = '  +10.50%  '
$ws.Range("E39").Value = '  +3.34%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.743'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9630'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.037'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +13.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '106.85'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.68%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000127'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.25%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.784'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.84%  '
$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D47").Value = '2.015.29'
$ws.Range("E47").Value = '  +1.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '65.38'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.5206'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.76%  '
$ws.Range("E50").Value = '  +4.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.036'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.71%  '
